# Refresh market-board derived columns (H:N) on the per-job Leve profit
# sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) with the latest pulled prices.
# One block per (sheet, leve row); values come from the scheduled price-
# refresh run. Cells with no applicable HQ/NQ price are cleared rather
# than left at 0, matching how the rest of the sheet is populated.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Cells.Item(33, 8).Value = 2337986.2  # H33 (currentAveragePrice)
$ws.Cells.Item(33, 9).Value = 3353362.2  # I33 (currentAveragePriceNQ)
$ws.Cells.Item(33, 11).Value = 3353362.2  # K33 (LevePriceNQ)
$ws.Cells.Item(33, 13).Value = -3353133.2  # M33 (LeveProfitNQ)

# Row 40: Stuck in the Moment / Horn Glue
$ws.Cells.Item(40, 8).Value = 1799.8  # H40 (currentAveragePrice)
$ws.Cells.Item(40, 9).Value = 1499.5  # I40 (currentAveragePriceNQ)
$ws.Cells.Item(40, 11).Value = 1499.5  # K40 (LevePriceNQ)
$ws.Cells.Item(40, 13).Value = -1324.5  # M40 (LeveProfitNQ)

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Cells.Item(74, 8).Value = 3289.4  # H74 (currentAveragePrice)
$ws.Cells.Item(74, 9).Value = 3289.4  # I74 (currentAveragePriceNQ)
$ws.Cells.Item(74, 11).Value = 3289.4  # K74 (LevePriceNQ)
$ws.Cells.Item(74, 13).Value = -2353.4  # M74 (LeveProfitNQ)

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Cells.Item(77, 8).Value = 3289.4  # H77 (currentAveragePrice)
$ws.Cells.Item(77, 9).Value = 3289.4  # I77 (currentAveragePriceNQ)
$ws.Cells.Item(77, 11).Value = 16447  # K77 (LevePriceNQ)
$ws.Cells.Item(77, 13).Value = -11767  # M77 (LeveProfitNQ)

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 5094.6416  # H138 (currentAveragePrice)
$ws.Cells.Item(138, 9).Value = 8347.294  # I138 (currentAveragePriceNQ)
$ws.Cells.Item(138, 10).Value = 3988.74  # J138 (currentAveragePriceHQ)
$ws.Cells.Item(138, 11).Value = 25041.882  # K138 (LevePriceNQ)
$ws.Cells.Item(138, 12).Value = 11966.22  # L138 (LevePriceHQ)
$ws.Cells.Item(138, 13).Value = -19901.882  # M138 (LeveProfitNQ)
$ws.Cells.Item(138, 14).Value = -22246.22  # N138 (LeveProfitHQ)

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Cells.Item(2, 8).Value = 1591.5555  # H2 (currentAveragePrice)
$ws.Cells.Item(2, 10).Value = 3691.6  # J2 (currentAveragePriceHQ)
$ws.Cells.Item(2, 12).Value = 3691.6  # L2 (LevePriceHQ)
$ws.Cells.Item(2, 14).Value = -3917.6  # N2 (LeveProfitHQ)

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 2357314.2  # H61 (currentAveragePrice)
$ws.Cells.Item(61, 9).Value = 4809  # I61 (currentAveragePriceNQ)
$ws.Cells.Item(61, 10).Value = 11767335  # J61 (currentAveragePriceHQ)
$ws.Cells.Item(61, 11).Value = 4809  # K61 (LevePriceNQ)
$ws.Cells.Item(61, 12).Value = 11767335  # L61 (LevePriceHQ)
$ws.Cells.Item(61, 13).Value = -4597  # M61 (LeveProfitNQ)
$ws.Cells.Item(61, 14).Value = -11767759  # N61 (LeveProfitHQ)

# Row 116: No Scope / Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 1591.5555  # H116 (currentAveragePrice)
$ws.Cells.Item(116, 10).Value = 3691.6  # J116 (currentAveragePriceHQ)
$ws.Cells.Item(116, 12).Value = 3691.6  # L116 (LevePriceHQ)
$ws.Cells.Item(116, 14).Value = -8279.6  # N116 (LeveProfitHQ)

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Cells.Item(122, 8).Value = 1375.7333  # H122 (currentAveragePrice)
$ws.Cells.Item(122, 9).Value = 1134.4878  # I122 (currentAveragePriceNQ)
$ws.Cells.Item(122, 10).Value = 3848.5  # J122 (currentAveragePriceHQ)
$ws.Cells.Item(122, 11).Value = 3403.463400000001  # K122 (LevePriceNQ)
$ws.Cells.Item(122, 12).Value = 11545.5  # L122 (LevePriceHQ)
$ws.Cells.Item(122, 13).Value = -953.4634000000005  # M122 (LeveProfitNQ)
$ws.Cells.Item(122, 14).Value = -16445.5  # N122 (LeveProfitHQ)

# Row 124: Ace of Gloves / High Durium Gauntlets of Fending
$ws.Cells.Item(124, 8).Value = 12357.25  # H124 (currentAveragePrice)
$ws.Cells.Item(124, 10).Value = 12357.25  # J124 (currentAveragePriceHQ)
$ws.Cells.Item(124, 12).Value = 12357.25  # L124 (LevePriceHQ)
$ws.Cells.Item(124, 14).Value = -22177.25  # N124 (LeveProfitHQ)

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 1717.6216  # H132 (currentAveragePrice)
$ws.Cells.Item(132, 9).Value = 1230.9656  # I132 (currentAveragePriceNQ)
$ws.Cells.Item(132, 10).Value = 3481.75  # J132 (currentAveragePriceHQ)
$ws.Cells.Item(132, 11).Value = 3692.8968  # K132 (LevePriceNQ)
$ws.Cells.Item(132, 12).Value = 10445.25  # L132 (LevePriceHQ)
$ws.Cells.Item(132, 13).Value = -1162.8968  # M132 (LeveProfitNQ)
$ws.Cells.Item(132, 14).Value = -15505.25  # N132 (LeveProfitHQ)

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 2357314.2  # H136 (currentAveragePrice)
$ws.Cells.Item(136, 9).Value = 4809  # I136 (currentAveragePriceNQ)
$ws.Cells.Item(136, 10).Value = 11767335  # J136 (currentAveragePriceHQ)
$ws.Cells.Item(136, 11).Value = 14427  # K136 (LevePriceNQ)
$ws.Cells.Item(136, 12).Value = 35302005  # L136 (LevePriceHQ)
$ws.Cells.Item(136, 13).Value = -11877  # M136 (LeveProfitNQ)
$ws.Cells.Item(136, 14).Value = -35307105  # N136 (LeveProfitHQ)

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Cells.Item(3, 8).Value = 1591.5555  # H3 (currentAveragePrice)
$ws.Cells.Item(3, 10).Value = 3691.6  # J3 (currentAveragePriceHQ)
$ws.Cells.Item(3, 12).Value = 3691.6  # L3 (LevePriceHQ)
$ws.Cells.Item(3, 14).Value = -3919.6  # N3 (LeveProfitHQ)

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Cells.Item(105, 8).Value = 5394.75  # H105 (currentAveragePrice)
$ws.Cells.Item(105, 9).Value = 5786.476  # I105 (currentAveragePriceNQ)
$ws.Cells.Item(105, 10).Value = 4219.5713  # J105 (currentAveragePriceHQ)
$ws.Cells.Item(105, 11).Value = 5786.476  # K105 (LevePriceNQ)
$ws.Cells.Item(105, 12).Value = 4219.5713  # L105 (LevePriceHQ)
$ws.Cells.Item(105, 13).Value = -4039.476  # M105 (LeveProfitNQ)
$ws.Cells.Item(105, 14).Value = -7713.5713  # N105 (LeveProfitHQ)

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 2800.6316  # H31 (currentAveragePrice)
$ws.Cells.Item(31, 9).Value = 2713.5186  # I31 (currentAveragePriceNQ)
$ws.Cells.Item(31, 11).Value = 2713.5186  # K31 (LevePriceNQ)
$ws.Cells.Item(31, 13).Value = -2418.5186  # M31 (LeveProfitNQ)

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 2800.6316  # H34 (currentAveragePrice)
$ws.Cells.Item(34, 9).Value = 2713.5186  # I34 (currentAveragePriceNQ)
$ws.Cells.Item(34, 11).Value = 2713.5186  # K34 (LevePriceNQ)
$ws.Cells.Item(34, 13).Value = -2511.5186  # M34 (LeveProfitNQ)

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Cells.Item(105, 8).Value = 3324.25  # H105 (currentAveragePrice)
$ws.Cells.Item(105, 9).Value = 2713  # I105 (currentAveragePriceNQ)
$ws.Cells.Item(105, 11).Value = 2713  # K105 (LevePriceNQ)
$ws.Cells.Item(105, 13).Value = -966  # M105 (LeveProfitNQ)

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Cells.Item(122, 8).Value = 1539.8  # H122 (currentAveragePrice)
$ws.Cells.Item(122, 10).Value = 826.5714  # J122 (currentAveragePriceHQ)
$ws.Cells.Item(122, 12).Value = 2479.7142  # L122 (LevePriceHQ)
$ws.Cells.Item(122, 14).Value = -7379.7142  # N122 (LeveProfitHQ)

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 19947.643  # H132 (currentAveragePrice)
$ws.Cells.Item(132, 9).Value = 25245.209  # I132 (currentAveragePriceNQ)
$ws.Cells.Item(132, 11).Value = 75735.62699999999  # K132 (LevePriceNQ)
$ws.Cells.Item(132, 13).Value = -73205.62699999999  # M132 (LeveProfitNQ)

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food / Table Salt
$ws.Cells.Item(2, 8).Value = 39.666668  # H2 (currentAveragePrice)
$ws.Cells.Item(2, 9).Value = 41.384617  # I2 (currentAveragePriceNQ)
$ws.Cells.Item(2, 11).Value = 248.307702  # K2 (LevePriceNQ)
$ws.Cells.Item(2, 13).Value = -135.307702  # M2 (LeveProfitNQ)

# Row 11: Putting the Squeeze On / Orange Juice
$ws.Cells.Item(11, 8).Value = 776.85  # H11 (currentAveragePrice)
$ws.Cells.Item(11, 9).Value = 681.2143  # I11 (currentAveragePriceNQ)
$ws.Cells.Item(11, 11).Value = 2043.6429  # K11 (LevePriceNQ)
$ws.Cells.Item(11, 13).Value = -1903.6429  # M11 (LeveProfitNQ)

# Row 21: Shy Is the Oyster / Raw Oyster
$ws.Cells.Item(21, 8).Value = 186  # H21 (currentAveragePrice)
$ws.Cells.Item(21, 9).Value = 186  # I21 (currentAveragePriceNQ)
$ws.Cells.Item(21, 11).Value = 558  # K21 (LevePriceNQ)
$ws.Cells.Item(21, 13).Value = -385  # M21 (LeveProfitNQ)

# Row 26: A Grape Idea / Grape Juice
$ws.Cells.Item(26, 9).Value = 530  # I26 (currentAveragePriceNQ)
$ws.Cells.Item(26, 10).Value = 421.33334  # J26 (currentAveragePriceHQ)
$ws.Cells.Item(26, 11).Value = 1590  # K26 (LevePriceNQ)
$ws.Cells.Item(26, 12).Value = 1264.00002  # L26 (LevePriceHQ)
$ws.Cells.Item(26, 13).Value = -1302  # M26 (LeveProfitNQ)
$ws.Cells.Item(26, 14).Value = -1840.00002  # N26 (LeveProfitHQ)

# Row 97: The Frier Never Lies / Cottonseed Oil
$ws.Cells.Item(97, 8).Value = 349.83334  # H97 (currentAveragePrice)
$ws.Cells.Item(97, 10).Value = 710  # J97 (currentAveragePriceHQ)
$ws.Cells.Item(97, 12).Value = 2130  # L97 (LevePriceHQ)
$ws.Cells.Item(97, 14).Value = -3122  # N97 (LeveProfitHQ)

# Row 98: Sweet Kiss of Death / Rice Vinegar
$ws.Cells.Item(98, 8).Value = 328.625  # H98 (currentAveragePrice)
$ws.Cells.Item(98, 10).Value = 328.625  # J98 (currentAveragePriceHQ)
$ws.Cells.Item(98, 12).Value = 985.875  # L98 (LevePriceHQ)
$ws.Cells.Item(98, 14).Value = -3981.875  # N98 (LeveProfitHQ)

# Row 113: Can't Eat Just One / Night Vinegar
$ws.Cells.Item(113, 8).Value = 787.6667  # H113 (currentAveragePrice)
$ws.Cells.Item(113, 10).Value = 950  # J113 (currentAveragePriceHQ)
$ws.Cells.Item(113, 12).Value = 2850  # L113 (LevePriceHQ)
$ws.Cells.Item(113, 14).Value = -7190  # N113 (LeveProfitHQ)

# Row 114: One Last Meal / Mushroom Saute
$ws.Cells.Item(114, 8).Value = 2203.3333  # H114 (currentAveragePrice)
$ws.Cells.Item(114, 9).Value = 2009.9166  # I114 (currentAveragePriceNQ)
$ws.Cells.Item(114, 11).Value = 6029.7498  # K114 (LevePriceNQ)
$ws.Cells.Item(114, 13).Value = -2775.7498  # M114 (LeveProfitNQ)

# Row 117: A Good Omen / Peppered Popotoes
$ws.Cells.Item(117, 8).Value = 66668320  # H117 (currentAveragePrice)
$ws.Cells.Item(117, 9).Value = 2064  # I117 (currentAveragePriceNQ)
$ws.Cells.Item(117, 11).Value = 6192  # K117 (LevePriceNQ)
$ws.Cells.Item(117, 13).Value = -2750  # M117 (LeveProfitNQ)

# Row 122: Salt of the North / Northern Sea Salt
$ws.Cells.Item(122, 8).Value = 2756011  # H122 (currentAveragePrice)
$ws.Cells.Item(122, 9).Value = 3788586  # I122 (currentAveragePriceNQ)
$ws.Cells.Item(122, 10).Value = 2478  # J122 (currentAveragePriceHQ)
$ws.Cells.Item(122, 11).Value = 34097274  # K122 (LevePriceNQ)
$ws.Cells.Item(122, 12).Value = 22302  # L122 (LevePriceHQ)
$ws.Cells.Item(122, 13).Value = -34094824  # M122 (LeveProfitNQ)
$ws.Cells.Item(122, 14).Value = -27202  # N122 (LeveProfitHQ)

# Row 129: Comfort Food / Yakow Moussaka
$ws.Cells.Item(129, 8).Value = 4686.3335  # H129 (currentAveragePrice)
$ws.Cells.Item(129, 9).Value = 2073.5  # I129 (currentAveragePriceNQ)
$ws.Cells.Item(129, 10).Value = 5432.857  # J129 (currentAveragePriceHQ)
$ws.Cells.Item(129, 11).Value = 6220.5  # K129 (LevePriceNQ)
$ws.Cells.Item(129, 12).Value = 16298.571  # L129 (LevePriceHQ)
$ws.Cells.Item(129, 13).Value = -1220.5  # M129 (LeveProfitNQ)
$ws.Cells.Item(129, 14).Value = -26298.571  # N129 (LeveProfitHQ)

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 3789953.8  # H131 (currentAveragePrice)
$ws.Cells.Item(131, 9).Value = 6994208  # I131 (currentAveragePriceNQ)
$ws.Cells.Item(131, 10).Value = 3108.0908  # J131 (currentAveragePriceHQ)
$ws.Cells.Item(131, 11).Value = 20982624  # K131 (LevePriceNQ)
$ws.Cells.Item(131, 12).Value = 9324.2724  # L131 (LevePriceHQ)
$ws.Cells.Item(131, 13).Value = -20977584  # M131 (LeveProfitNQ)
$ws.Cells.Item(131, 14).Value = -19404.2724  # N131 (LeveProfitHQ)

# Row 132: More Mezcal / Cooking Mezcal
$ws.Cells.Item(132, 8).Value = 1999.7142  # H132 (currentAveragePrice)
$ws.Cells.Item(132, 9).Value = 0  # I132 (currentAveragePriceNQ)
$ws.Cells.Item(132, 10).Value = 1999.7142  # J132 (currentAveragePriceHQ)
$ws.Cells.Item(132, 11).Value = 0  # K132 (LevePriceNQ)
$ws.Cells.Item(132, 12).Value = 17997.4278  # L132 (LevePriceHQ)
$ws.Cells.Item(132, 13).ClearContents()  # M132 (LeveProfitNQ) no longer applicable
$ws.Cells.Item(132, 14).Value = -23057.4278  # N132 (LeveProfitHQ)

# Row 141: Ocean Explosion / Acqua Pazza
$ws.Cells.Item(141, 8).Value = 7619.6665  # H141 (currentAveragePrice)
$ws.Cells.Item(141, 9).Value = 7209.5  # I141 (currentAveragePriceNQ)
$ws.Cells.Item(141, 10).Value = 8440  # J141 (currentAveragePriceHQ)
$ws.Cells.Item(141, 11).Value = 21628.5  # K141 (LevePriceNQ)
$ws.Cells.Item(141, 12).Value = 25320  # L141 (LevePriceHQ)
$ws.Cells.Item(141, 13).Value = -16448.5  # M141 (LeveProfitNQ)
$ws.Cells.Item(141, 14).Value = -35680  # N141 (LeveProfitHQ)

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... / Koppranickel Ingot
$ws.Cells.Item(97, 8).Value = 517.64514  # H97 (currentAveragePrice)
$ws.Cells.Item(97, 9).Value = 480.24  # I97 (currentAveragePriceNQ)
$ws.Cells.Item(97, 11).Value = 480.24  # K97 (LevePriceNQ)
$ws.Cells.Item(97, 13).Value = 15.75999999999999  # M97 (LeveProfitNQ)

# Row 105: Untucked / Palladium Tuck
$ws.Cells.Item(105, 8).Value = 76950  # H105 (currentAveragePrice)
$ws.Cells.Item(105, 10).Value = 76950  # J105 (currentAveragePriceHQ)
$ws.Cells.Item(105, 12).Value = 76950  # L105 (LevePriceHQ)
$ws.Cells.Item(105, 14).Value = -83938  # N105 (LeveProfitHQ)

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 3087.375  # H126 (currentAveragePrice)
$ws.Cells.Item(126, 9).Value = 2580  # I126 (currentAveragePriceNQ)
$ws.Cells.Item(126, 10).Value = 3933  # J126 (currentAveragePriceHQ)
$ws.Cells.Item(126, 11).Value = 7740  # K126 (LevePriceNQ)
$ws.Cells.Item(126, 12).Value = 11799  # L126 (LevePriceHQ)
$ws.Cells.Item(126, 13).Value = -5270  # M126 (LeveProfitNQ)
$ws.Cells.Item(126, 14).Value = -16739  # N126 (LeveProfitHQ)

$ws = $wb.Worksheets.Item("LTW")
# Row 34: Breeches Served Cold / Goatskin Breeches
$ws.Cells.Item(34, 8).Value = 9000  # H34 (currentAveragePrice)
$ws.Cells.Item(34, 9).Value = 9000  # I34 (currentAveragePriceNQ)
$ws.Cells.Item(34, 11).Value = 9000  # K34 (LevePriceNQ)
$ws.Cells.Item(34, 13).Value = -8828  # M34 (LeveProfitNQ)

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Cells.Item(61, 8).Value = 2884.121  # H61 (currentAveragePrice)
$ws.Cells.Item(61, 9).Value = 2798.6191  # I61 (currentAveragePriceNQ)
$ws.Cells.Item(61, 11).Value = 2798.6191  # K61 (LevePriceNQ)
$ws.Cells.Item(61, 13).Value = -2596.6191  # M61 (LeveProfitNQ)

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Cells.Item(113, 8).Value = 2884.121  # H113 (currentAveragePrice)
$ws.Cells.Item(113, 9).Value = 2798.6191  # I113 (currentAveragePriceNQ)
$ws.Cells.Item(113, 11).Value = 2798.6191  # K113 (LevePriceNQ)
$ws.Cells.Item(113, 13).Value = -628.6190999999999  # M113 (LeveProfitNQ)

$ws = $wb.Worksheets.Item("WVR")
# Row 4: Not Cool Enough / Hempen Undershirt
$ws.Cells.Item(4, 8).Value = 25000.666  # H4 (currentAveragePrice)
$ws.Cells.Item(4, 9).Value = 30002  # I4 (currentAveragePriceNQ)
$ws.Cells.Item(4, 11).Value = 30002  # K4 (LevePriceNQ)
$ws.Cells.Item(4, 13).Value = -29889  # M4 (LeveProfitNQ)

# Row 34: He's Got Legs / Velveteen Sarouel
$ws.Cells.Item(34, 8).Value = 0  # H34 (currentAveragePrice)
$ws.Cells.Item(34, 9).Value = 0  # I34 (currentAveragePriceNQ)
$ws.Cells.Item(34, 11).Value = 0  # K34 (LevePriceNQ)
$ws.Cells.Item(34, 13).ClearContents()  # M34 (LeveProfitNQ) no longer applicable

# Row 45: Private Concerns / Linen Trousers
$ws.Cells.Item(45, 8).Value = 25744.25  # H45 (currentAveragePrice)
$ws.Cells.Item(45, 10).Value = 29988.5  # J45 (currentAveragePriceHQ)
$ws.Cells.Item(45, 12).Value = 29988.5  # L45 (LevePriceHQ)
$ws.Cells.Item(45, 14).Value = -30970.5  # N45 (LeveProfitHQ)

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 2129.889  # H132 (currentAveragePrice)
$ws.Cells.Item(132, 9).Value = 1184.4  # I132 (currentAveragePriceNQ)
$ws.Cells.Item(132, 10).Value = 4278.727  # J132 (currentAveragePriceHQ)
$ws.Cells.Item(132, 11).Value = 3553.2  # K132 (LevePriceNQ)
$ws.Cells.Item(132, 12).Value = 12836.181  # L132 (LevePriceHQ)
$ws.Cells.Item(132, 13).Value = -1023.2  # M132 (LeveProfitNQ)
$ws.Cells.Item(132, 14).Value = -17896.181  # N132 (LeveProfitHQ)

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 33645.613  # H136 (currentAveragePrice)
$ws.Cells.Item(136, 9).Value = 46616.5  # I136 (currentAveragePriceNQ)
$ws.Cells.Item(136, 11).Value = 139849.5  # K136 (LevePriceNQ)
$ws.Cells.Item(136, 13).Value = -137299.5  # M136 (LeveProfitNQ)
